$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header for new Date column
$ws.Range("C1").Value = "Date"

# Convert Price to numeric
$ws.Range("B2").Value = 54590

# Set Date value and format
$ws.Range("C2").NumberFormat = "yyyy-mm-dd"
$ws.Range("C2").Value = 45326
